$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update row 2 (Chaitanya Somawar)
$ws1.Range("C2").Value = "chaitanya234aa45@gmail.com"
$ws1.Range("D2").Value = "'334534112"
$ws1.Range("E2").Value = "test@322121"

# Update row 3 (Madhurima)
$ws1.Range("C3").Value = "madhu898jkhkj179@gmail.com"
$ws1.Range("D3").Value = "'452452323233"
$ws1.Range("E3").Value = "test@23412122"

# Update row 4 (Vaibhav)
$ws1.Range("C4").Value = "vaibhav9utiut6549732@gmail.com"
$ws1.Range("D4").Value = "'9879378711"

# Update row 5 (Nikhil)
$ws1.Range("C5").Value = "nikhilq3hjkytu7653445@test.com"
$ws1.Range("D5").Value = "'232387987"

# Add a new worksheet "LoginData" after the existing sheet
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "LoginData"

# Restore active sheet/selection on the first sheet
$ws1.Activate()
$ws1.Range("C5").Select()
